# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the refreshed counts captured at the new scrape time.

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> { row -> new F value }
$updates = @{
    "展览" = @{
        2  = 4339
        3  = 2446
        4  = 482
        5  = 23
        6  = 47
        7  = 59
        8  = 215
        10 = 148
        11 = 156
        12 = 1614
        13 = 296
        14 = 3401
        15 = 229
    }
    "全部类型" = @{
        2  = 4339
        3  = 2446
        4  = 482
        5  = 23
        7  = 47
        8  = 59
        10 = 215
        12 = 148
        13 = 156
        16 = 1614
        17 = 296
        18 = 3401
        19 = 229
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsMap = $updates[$sheetName]
    foreach ($row in $rowsMap.Keys) {
        $ws.Cells.Item($row, 6).Value = $rowsMap[$row]
    }
}
